$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "H 72" record (row 2) entirely, shifting all rows below it up by one.
$ws.Rows.Item(2).Delete()
